$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 20323
$ws1.Range("F6").Value = 1107
$ws1.Range("F8").Value = 7647
$ws1.Range("F10").Value = 739
$ws1.Range("F11").Value = 278
$ws1.Range("F13").Value = 164
$ws1.Range("F14").Value = 131
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 238
$ws1.Range("F18").Value = 1343
$ws1.Range("F19").Value = 451
$ws1.Range("F21").Value = 690
$ws1.Range("F24").Value = 72
$ws1.Range("F25").Value = 329
$ws1.Range("F29").Value = 190
$ws1.Range("F32").Value = 83
$ws1.Range("F33").Value = 3342
$ws1.Range("F35").Value = 88
$ws1.Range("F36").Value = 53
$ws1.Range("F37").Value = 12730
$ws1.Range("F39").Value = 93
$ws1.Range("F40").Value = 33
$ws1.Range("F43").Value = 383
$ws1.Range("F44").Value = 4011

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 206

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 20323
$ws4.Range("F6").Value = 1107
$ws4.Range("F8").Value = 7647
$ws4.Range("F10").Value = 739
$ws4.Range("F11").Value = 278
$ws4.Range("F13").Value = 164
$ws4.Range("F14").Value = 131
$ws4.Range("F15").Value = 17
$ws4.Range("F16").Value = 238
$ws4.Range("F18").Value = 1343
$ws4.Range("F19").Value = 451
$ws4.Range("F21").Value = 690
$ws4.Range("F24").Value = 72
$ws4.Range("F25").Value = 329
$ws4.Range("F29").Value = 190
$ws4.Range("F30").Value = 206
$ws4.Range("F34").Value = 83
$ws4.Range("F36").Value = 3345
$ws4.Range("F38").Value = 88
$ws4.Range("F39").Value = 53
$ws4.Range("F40").Value = 12730
$ws4.Range("F42").Value = 93
$ws4.Range("F43").Value = 33
$ws4.Range("F46").Value = 383
$ws4.Range("F47").Value = 4011
